$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 33 (item 5512)
$ws.Cells.Item(33, 8).Value = 177.6
$ws.Cells.Item(33, 9).Value = 177.6
$ws.Cells.Item(33, 11).Value = 177.6
$ws.Cells.Item(33, 13).Value = 51.40000000000001

# Row 41 (item 5478)
$ws.Cells.Item(41, 8).Value = 399.92856
$ws.Cells.Item(41, 9).Value = 244.77777
$ws.Cells.Item(41, 10).Value = 679.2
$ws.Cells.Item(41, 11).Value = 244.77777
$ws.Cells.Item(41, 12).Value = 679.2
$ws.Cells.Item(41, 13).Value = 195.22223
$ws.Cells.Item(41, 14).Value = -1559.2

# Row 53 (item 5479)
$ws.Cells.Item(53, 8).Value = 199.3158
$ws.Cells.Item(53, 9).Value = 203.33333
$ws.Cells.Item(53, 10).Value = 192.42857
$ws.Cells.Item(53, 11).Value = 203.33333
$ws.Cells.Item(53, 12).Value = 192.42857
$ws.Cells.Item(53, 13).Value = 433.66667
$ws.Cells.Item(53, 14).Value = -1466.42857

# Row 129 (item 36115)
$ws.Cells.Item(129, 8).Value = 4117.857

# Row 132 (item 44049)
$ws.Cells.Item(132, 8).Value = 11804.053
$ws.Cells.Item(132, 9).Value = 11804.053
$ws.Cells.Item(132, 10).Value = 0
$ws.Cells.Item(132, 11).Value = 35412.159
$ws.Cells.Item(132, 12).Value = 0
$ws.Cells.Item(132, 13).Value = -32882.159
$ws.Cells.Item(132, 14).ClearContents()

# Row 135 (item 44047)
$ws.Cells.Item(135, 8).Value = 2120.3333
$ws.Cells.Item(135, 9).Value = 2120.3333
$ws.Cells.Item(135, 11).Value = 19082.9997
$ws.Cells.Item(135, 13).Value = -16547.9997

# Row 137 (item 44013)
$ws.Cells.Item(137, 8).Value = 8320.429
$ws.Cells.Item(137, 9).Value = 1323.3334
$ws.Cells.Item(137, 11).Value = 3970.0002
$ws.Cells.Item(137, 13).Value = -1420.0002

# Row 138 (item 44169)
$ws.Cells.Item(138, 8).Value = 3228.0645
$ws.Cells.Item(138, 10).Value = 2396.625
$ws.Cells.Item(138, 12).Value = 7189.875
$ws.Cells.Item(138, 14).Value = -17469.875

$ws = $wb.Worksheets.Item("ARM")
# Row 61 (item 43999)
$ws.Cells.Item(61, 8).Value = 2269.111
$ws.Cells.Item(61, 9).Value = 2014.5238
$ws.Cells.Item(61, 10).Value = 5833.3335
$ws.Cells.Item(61, 11).Value = 2014.5238
$ws.Cells.Item(61, 12).Value = 5833.3335
$ws.Cells.Item(61, 13).Value = -1802.5238
$ws.Cells.Item(61, 14).Value = -6257.3335

# Row 122 (item 36168)
$ws.Cells.Item(122, 8).Value = 2122.0908
$ws.Cells.Item(122, 9).Value = 1984.3
$ws.Cells.Item(122, 10).Value = 3500
$ws.Cells.Item(122, 11).Value = 5952.9
$ws.Cells.Item(122, 12).Value = 10500
$ws.Cells.Item(122, 13).Value = -3502.9
$ws.Cells.Item(122, 14).Value = -15400

# Row 136 (item 43999)
$ws.Cells.Item(136, 8).Value = 2269.111
$ws.Cells.Item(136, 9).Value = 2014.5238
$ws.Cells.Item(136, 10).Value = 5833.3335
$ws.Cells.Item(136, 11).Value = 6043.5714
$ws.Cells.Item(136, 12).Value = 17500.0005
$ws.Cells.Item(136, 13).Value = -3493.5714
$ws.Cells.Item(136, 14).Value = -22600.0005

$ws = $wb.Worksheets.Item("BSM")
# Row 6 (item 27119)
$ws.Cells.Item(6, 8).Value = 36852.57
$ws.Cells.Item(6, 10).Value = 36852.57
$ws.Cells.Item(6, 12).Value = 36852.57
$ws.Cells.Item(6, 14).Value = -37078.57

# Row 134 (item 43998)
$ws.Cells.Item(134, 8).Value = 2574.258
$ws.Cells.Item(134, 9).Value = 2574.258
$ws.Cells.Item(134, 11).Value = 7722.773999999999
$ws.Cells.Item(134, 13).Value = -5187.773999999999

$ws = $wb.Worksheets.Item("CRP")
# Row 22 (item 5367)
$ws.Cells.Item(22, 8).Value = 1858.3334
$ws.Cells.Item(22, 9).Value = 821.36365
$ws.Cells.Item(22, 11).Value = 821.36365
$ws.Cells.Item(22, 13).Value = -471.36365

# Row 31 (item 44023)
$ws.Cells.Item(31, 8).Value = 2557.3845
$ws.Cells.Item(31, 9).Value = 2746.6
$ws.Cells.Item(31, 10).Value = 1926.6666
$ws.Cells.Item(31, 11).Value = 2746.6
$ws.Cells.Item(31, 12).Value = 1926.6666
$ws.Cells.Item(31, 13).Value = -2451.6
$ws.Cells.Item(31, 14).Value = -2516.6666

# Row 34 (item 44023)
$ws.Cells.Item(34, 8).Value = 2557.3845
$ws.Cells.Item(34, 9).Value = 2746.6
$ws.Cells.Item(34, 10).Value = 1926.6666
$ws.Cells.Item(34, 11).Value = 2746.6
$ws.Cells.Item(34, 12).Value = 1926.6666
$ws.Cells.Item(34, 13).Value = -2544.6
$ws.Cells.Item(34, 14).Value = -2330.6666

# Row 37 (item 2021)
$ws.Cells.Item(37, 8).Value = 3000
$ws.Cells.Item(37, 9).Value = 3000
$ws.Cells.Item(37, 11).Value = 3000
$ws.Cells.Item(37, 13).Value = -2893

# Row 58 (item 44021)
$ws.Cells.Item(58, 8).Value = 3273.4333
$ws.Cells.Item(58, 9).Value = 2738.3333
$ws.Cells.Item(58, 10).Value = 4522
$ws.Cells.Item(58, 11).Value = 2738.3333
$ws.Cells.Item(58, 12).Value = 4522
$ws.Cells.Item(58, 13).Value = -2535.3333
$ws.Cells.Item(58, 14).Value = -4928

# Row 105 (item 19928)
$ws.Cells.Item(105, 8).Value = 10107.728
$ws.Cells.Item(105, 9).Value = 12220.889
$ws.Cells.Item(105, 10).Value = 598.5
$ws.Cells.Item(105, 11).Value = 12220.889
$ws.Cells.Item(105, 12).Value = 598.5
$ws.Cells.Item(105, 13).Value = -10473.889
$ws.Cells.Item(105, 14).Value = -4092.5

# Row 136 (item 44021)
$ws.Cells.Item(136, 8).Value = 3273.4333
$ws.Cells.Item(136, 9).Value = 2738.3333
$ws.Cells.Item(136, 10).Value = 4522
$ws.Cells.Item(136, 11).Value = 8214.999899999999
$ws.Cells.Item(136, 12).Value = 13566
$ws.Cells.Item(136, 13).Value = -5664.999899999999
$ws.Cells.Item(136, 14).Value = -18666

$ws = $wb.Worksheets.Item("CUL")
# Row 69 (item 12850)
$ws.Cells.Item(69, 8).Value = 4112.143
$ws.Cells.Item(69, 10).Value = 4900
$ws.Cells.Item(69, 12).Value = 14700
$ws.Cells.Item(69, 14).Value = -16322

# Row 72 (item 12850)
$ws.Cells.Item(72, 8).Value = 4112.143
$ws.Cells.Item(72, 10).Value = 4900
$ws.Cells.Item(72, 12).Value = 44100
$ws.Cells.Item(72, 14).Value = -52212

# Row 75 (item 12863)
$ws.Cells.Item(75, 8).Value = 497.5
$ws.Cells.Item(75, 9).Value = 497.5
$ws.Cells.Item(75, 11).Value = 1492.5
$ws.Cells.Item(75, 13).Value = -494.5

# Row 78 (item 12863)
$ws.Cells.Item(78, 8).Value = 497.5
$ws.Cells.Item(78, 9).Value = 497.5
$ws.Cells.Item(78, 11).Value = 4477.5
$ws.Cells.Item(78, 13).Value = 514.5

# Row 108 (item 27853)
$ws.Cells.Item(108, 8).Value = 8000
$ws.Cells.Item(108, 9).Value = 0
$ws.Cells.Item(108, 10).Value = 8000
$ws.Cells.Item(108, 11).Value = 0
$ws.Cells.Item(108, 12).Value = 24000
$ws.Cells.Item(108, 13).ClearContents()
$ws.Cells.Item(108, 14).Value = -29760

# Row 121 (item 27878)
$ws.Cells.Item(121, 8).Value = 17635.812
$ws.Cells.Item(121, 9).Value = 240.42857
$ws.Cells.Item(121, 10).Value = 31165.555
$ws.Cells.Item(121, 11).Value = 721.28571
$ws.Cells.Item(121, 12).Value = 93496.66500000001
$ws.Cells.Item(121, 13).Value = 588.71429
$ws.Cells.Item(121, 14).Value = -96116.66500000001

# Row 137 (item 44088)
$ws.Cells.Item(137, 8).Value = 8668.429
$ws.Cells.Item(137, 9).Value = 2622.5715
$ws.Cells.Item(137, 10).Value = 14714.286
$ws.Cells.Item(137, 11).Value = 7867.7145
$ws.Cells.Item(137, 12).Value = 44142.858
$ws.Cells.Item(137, 13).Value = -2767.7145
$ws.Cells.Item(137, 14).Value = -54342.858

$ws = $wb.Worksheets.Item("GSM")
# Row 15 (item 12018)
$ws.Cells.Item(15, 8).Value = 22039.25
$ws.Cells.Item(15, 10).Value = 19473.428
$ws.Cells.Item(15, 12).Value = 19473.428
$ws.Cells.Item(15, 14).Value = -20049.428

# Row 18 (item 4309)
$ws.Cells.Item(18, 8).Value = 7332.6665
$ws.Cells.Item(18, 10).Value = 7332.6665
$ws.Cells.Item(18, 12).Value = 7332.6665
$ws.Cells.Item(18, 14).Value = -7918.6665

# Row 81 (item 12018)
$ws.Cells.Item(81, 8).Value = 22039.25
$ws.Cells.Item(81, 10).Value = 19473.428
$ws.Cells.Item(81, 12).Value = 19473.428
$ws.Cells.Item(81, 14).Value = -21469.428

# Row 84 (item 12018)
$ws.Cells.Item(84, 8).Value = 22039.25
$ws.Cells.Item(84, 10).Value = 19473.428
$ws.Cells.Item(84, 12).Value = 58420.284
$ws.Cells.Item(84, 14).Value = -68404.284

# Row 132 (item 44008)
$ws.Cells.Item(132, 8).Value = 21456
$ws.Cells.Item(132, 9).Value = 21456
$ws.Cells.Item(132, 10).Value = 0
$ws.Cells.Item(132, 11).Value = 64368
$ws.Cells.Item(132, 12).Value = 0
$ws.Cells.Item(132, 13).Value = -61838
$ws.Cells.Item(132, 14).ClearContents()

$ws = $wb.Worksheets.Item("LTW")
# Row 40 (item 36248)
$ws.Cells.Item(40, 8).Value = 4471.8
$ws.Cells.Item(40, 9).Value = 3853.353
$ws.Cells.Item(40, 10).Value = 7976.3335
$ws.Cells.Item(40, 11).Value = 3853.353
$ws.Cells.Item(40, 12).Value = 7976.3335
$ws.Cells.Item(40, 13).Value = -3717.353
$ws.Cells.Item(40, 14).Value = -8248.333500000001

# Row 136 (item 44060)
$ws.Cells.Item(136, 8).Value = 2951.6924
$ws.Cells.Item(136, 9).Value = 2689.4348
$ws.Cells.Item(136, 11).Value = 8068.3044
$ws.Cells.Item(136, 13).Value = -5518.3044

$ws = $wb.Worksheets.Item("WVR")
# Row 136 (item 44031)
$ws.Cells.Item(136, 8).Value = 2492.6
$ws.Cells.Item(136, 9).Value = 2136.9092
$ws.Cells.Item(136, 11).Value = 6410.7276
$ws.Cells.Item(136, 13).Value = -3860.7276
